$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect before editing, then restore protection after.
$ws.Unprotect("D382")

# Update the confidential disclaimer text (date changed from 2021-05-19 to 2021-05-20)
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-20 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-39
$ws.Range("D2").Value = 0.05724175142782228
$ws.Range("E2").Value = 0.0210121100328815
$ws.Range("D3").Value = 0.05191149252800769
$ws.Range("E3").Value = 0.01382033563672258
$ws.Range("D4").Value = 0.317987784827803
$ws.Range("E4").Value = 0.009870740305522929
$ws.Range("D5").Value = 0.03450303585719301
$ws.Range("E5").Value = 0.004913670400396031
$ws.Range("D6").Value = 0.03117070323076017
$ws.Range("E6").Value = 0.008237232289950436
$ws.Range("D7").Value = 0.0309604907250271
$ws.Range("E7").Value = -0.001737943020296839
$ws.Range("D8").Value = 0.02905267077710941
$ws.Range("E8").Value = 0.005820790216368721
$ws.Range("D9").Value = 0.02424244494759741
$ws.Range("E9").Value = 0.003523111612175889
$ws.Range("D10").Value = 0.02425077230942939
$ws.Range("E10").Value = 0.01560642747083407
$ws.Range("D11").Value = 0.02343543817723871
$ws.Range("E11").Value = 0.01600816352562284
$ws.Range("D12").Value = 0.02329995840589526
$ws.Range("E12").Value = -0.002382654276864571
$ws.Range("D13").Value = 0.01987858706448965
$ws.Range("E13").Value = 0.01234713770898566
$ws.Range("D14").Value = 0.02187533896633098
$ws.Range("E14").Value = 0.01024890190336758
$ws.Range("D15").Value = 0.02014260713898324
$ws.Range("E15").Value = 0.004611225950071685
$ws.Range("D16").Value = 0.02185932480896178
$ws.Range("E16").Value = -0.002100122100122226
$ws.Range("D17").Value = 0.01926930175711605
$ws.Range("E17").Value = 0.01861598980552936
$ws.Range("D18").Value = 0.01418470003134504
$ws.Range("E18").Value = 0.01065751445086738
$ws.Range("D19").Value = 0.01621401405317042
$ws.Range("E19").Value = 0.01474926253687303
$ws.Range("D20").Value = 0.01519209729091703
$ws.Range("E20").Value = -0.002283907238229133
$ws.Range("D21").Value = 0.01637159336168337
$ws.Range("E21").Value = -0.002373685995252561
$ws.Range("D22").Value = 0.01203111614833466
$ws.Range("E22").Value = 0.04138714371916352
$ws.Range("D23").Value = 0.01503643968128838
$ws.Range("E23").Value = 0.008860993169651055
$ws.Range("D24").Value = 0.01329473992581388
$ws.Range("E24").Value = 0.02348066298342544
$ws.Range("D25").Value = 0.01400438061936782
$ws.Range("E25").Value = 0.007242233657327901
$ws.Range("D26").Value = 0.01366733598727067
$ws.Range("E26").Value = 0.01220140917683454
$ws.Range("D27").Value = 0.0127754541828552
$ws.Range("E27").Value = 0.002331528279181594
$ws.Range("D28").Value = 0.01322128832401381
$ws.Range("E28").Value = 0.003488372093023218
$ws.Range("D29").Value = 0.01445779479501452
$ws.Range("E29").Value = 0.007280944012051149
$ws.Range("D30").Value = 0.01345018401334428
$ws.Range("E30").Value = -0.0002540005080009511
$ws.Range("D31").Value = 0.01242100749975018
$ws.Range("E31").Value = 0.0123770886337069
$ws.Range("D32").Value = 0.01360509429562903
$ws.Range("E32").Value = 0.01087613293051359
$ws.Range("D33").Value = 0.01270883528819932
$ws.Range("E33").Value = 0.003780241935483764
$ws.Range("D34").Value = 0.006006696907089703
$ws.Range("E34").Value = 0.03887101647619229
$ws.Range("D35").Value = 0.00520673636597346
$ws.Range("E35").Value = 0.02864465860159937
$ws.Range("D36").Value = 0.005223391089637431
$ws.Range("E36").Value = 0.02824673997465554
$ws.Range("D37").Value = 0.005129548127453904
$ws.Range("E37").Value = 0.02331050846046567
$ws.Range("D38").Value = 0.004715849062082831
$ws.Range("E38").Value = 0.03006429412297384
$ws.Range("E39").Value = 0.01010674823780211

# Restore sheet protection
$ws.Protect("D382", $true, $true, $true)

